$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Capture text that must survive from rows that will be removed ---
# (Read via .Value2 and re-write via .Value so the underlying shared-string
#  bytes are copied verbatim - avoids any console/string-literal transcoding
#  of the accented characters.)
$docJuan    = $ws.Range("C24").Value2   # "9100813"
$nameJuan   = $ws.Range("D24").Value2   # "JUAN CARLOS SILVA BOLA...OS"
$docJunior  = $ws.Range("C26").Value2   # "1002244636"
$nameJunior = $ws.Range("D26").Value2   # "JUNIOR JOSE CANTILLO TORRES"

# --- 2. Give row 18 the heavier "last row of table" border/style that
#        currently belongs to row 30 (the row it is about to replace). ---
$ws.Rows("30:30").Copy()
$ws.Rows("18:18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Update the worker/contribution table (rows 16-18) ---
$ws.Range("G16").Value = 781243

$ws.Range("C17").Value = $docJuan
$ws.Range("D17").Value = $nameJuan
$ws.Range("E17").Value = "2001"
$ws.Range("F17").Value = 32021
$ws.Range("G17").Value = 828116

$ws.Range("C18").Value = $docJunior
$ws.Range("D18").Value = $nameJunior
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# --- 4. Update the summary fields above the table ---
$ws.Range("E11").Value = 91044
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 2

# --- 5. Drop the now-unused worker rows (19-30); this shifts the
#        signature block (old rows 35/36) up to rows 23/24. ---
$ws.Range("19:30").EntireRow.Delete()

# --- 6. Column D is now narrower since the longest name got shorter ---
$ws.Columns("D:D").ColumnWidth = 33.36328125
